$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit order matters for shared-string table ordering: Careers+ (C18) first,
# then Company (C16), then News&Update (C15) -- matching how the source
# workbook's shared strings were appended.

# Update the "Careers+" expected-result cell (row 18, column C)
$ws.Range("C18").Value = "แสดงข้อมูลในหน้า Careers+ โดยมีข้อมูลดังนี้`n- Life at Swift Dynamics`n    - แสดงรูปภาพ , ชื่อ และตำแหน่งการทำงานของแต่ละท่าน`n- LOOKING FOR THE RIGHT POSITION`n    - เมื่อกดเข้าไปจะแสดงตำแหน่งงานที่เปิดรับสมัคร`n"

# Update the "Company" expected-result cell (row 16, column C)
$ws.Range("C16").Value = "แสดงข้อมูลในหน้า Company โดยมีข้อมูลดังนี้`n- Why Swift Dynamics?`n   - แสดงรายละเอียดข้อมูล`n- แสดงรูปภาพ คำคม ชื่อและตำแหน่งงาน`n- Maps , Address , Tel , Business Hours`n   - แสดงรูปภาพแผนที่ ที่อยู่แผนที่ เบอร์โทรติดต่อ และเวลาเปิด-ปิด ทำการ"

# Update the "News&Update" expected-result cell (row 15, column C)
$ws.Range("C15").Value = "ระบบแสดงข้อมูลรายการ News&Update ล่าสุด`n- แสดงรูปภาพ หัวข้อข่าว วันที่ และ รายละเอียดเพิ่มเติมของข่าว"
$ws.Range("C15").WrapText = $true

# Scroll the view and select C15, matching the final state of the sheet view
$excel.ActiveWindow.ScrollRow = 15
$ws.Range("C15").Select()
